$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 95; existing rows 95-179 shift down to 96-180.
$ws.Rows("95").Insert()

# Populate the newly inserted row 95 with the new weekly price record.
$ws.Range("A95").Value = 7
$ws.Range("B95").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C95").Value = "Ñuble"
$ws.Range("D95").Value = 44512
$ws.Range("E95").Value = 16
$ws.Range("F95").Value = 100112009
$ws.Range("G95").Value = "Acelga"
$ws.Range("H95").Value = "Sin especificar"
$ws.Range("I95").Value = "Primera"
$ws.Range("J95").Value = 400
$ws.Range("K95").Value = 350
$ws.Range("L95").Value = 400
$ws.Range("M95").Value = 375
$ws.Range("N95").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O95").Value = "Provincia de Diguillín"
$ws.Range("P95").Value = 375
$ws.Range("Q95").Value = 1
$ws.Range("R95").Value = "Hortaliza"
